$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1831553333333333
$ws.Range("H2").Value = 0.549466
$ws.Range("I2").Value = 0.3807808465430998
$ws.Range("J2").Value = 0.3807808465430998
$ws.Range("M2").Value = 5.560959666666666
$ws.Range("N2").Value = 16.682879
$ws.Range("O2").Value = 0.1423184612405788
$ws.Range("P2").Value = 0.1423184612405788
$ws.Range("Q2").Value = 1.018519421401556
$ws.Range("R2").Value = 9.166674792614
$ws.Range("S2").Value = 0.05419214414989895
$ws.Range("T2").Value = 0.05419214414989892
# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1831553333333333
$ws.Range("H3").Value = 0.549466
$ws.Range("I3").Value = 0.3807808465430998
$ws.Range("J3").Value = 0.3807808465430998
$ws.Range("M3").Value = 5.635187000000001
$ws.Range("O3").Value = 0.1442181189427041
$ws.Range("P3").Value = 0.1442181189427041
$ws.Range("Q3").Value = 1.032114553380667
$ws.Range("R3").Value = 9.289030980426002
$ws.Range("S3").Value = 0.05491549741785635
$ws.Range("T3").Value = 0.05491549741785633
# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1831553333333333
$ws.Range("H4").Value = 0.549466
$ws.Range("I4").Value = 0.3807808465430998
$ws.Range("J4").Value = 0.3807808465430998
$ws.Range("M4").Value = 7.323012666666667
$ws.Range("N4").Value = 21.969038
$ws.Range("O4").Value = 0.1874136762063552
$ws.Range("P4").Value = 0.1874136762063552
$ws.Range("Q4").Value = 1.341248825967556
$ws.Range("R4").Value = 12.071239433708
$ws.Range("S4").Value = 0.07136353827961035
$ws.Range("T4").Value = 0.07136353827961034
# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1831553333333333
$ws.Range("H5").Value = 0.549466
$ws.Range("I5").Value = 0.3807808465430998
$ws.Range("J5").Value = 0.3807808465430998
$ws.Range("M5").Value = 20.55489766666667
$ws.Range("N5").Value = 61.664693
$ws.Range("O5").Value = 0.5260497436103619
$ws.Range("P5").Value = 0.5260497436103618
$ws.Range("Q5").Value = 3.764739133770889
$ws.Range("R5").Value = 33.882652203938
$ws.Range("S5").Value = 0.2003096666957342
$ws.Range("T5").Value = 0.2003096666957342
# Row 6
$ws.Range("I6").Value = 0.3553476858595785
$ws.Range("J6").Value = 0.3553476858595784
$ws.Range("M6").Value = 5.560959666666666
$ws.Range("N6").Value = 16.682879
$ws.Range("O6").Value = 0.1423184612405788
$ws.Range("P6").Value = 0.1423184612405788
$ws.Range("Q6").Value = 0.9504903481460001
$ws.Range("R6").Value = 8.554413133314
$ws.Range("S6").Value = 0.05057253585693579
$ws.Range("T6").Value = 0.05057253585693578
# Row 7
$ws.Range("I7").Value = 0.3553476858595785
$ws.Range("J7").Value = 0.3553476858595784
$ws.Range("M7").Value = 5.635187000000001
$ws.Range("O7").Value = 0.1442181189427041
$ws.Range("P7").Value = 0.1442181189427041
$ws.Range("Q7").Value = 0.9631774324140003
$ws.Range("R7").Value = 8.668596891726002
$ws.Range("S7").Value = 0.05124757482531136
$ws.Range("T7").Value = 0.05124757482531134
# Row 8
$ws.Range("I8").Value = 0.3553476858595785
$ws.Range("J8").Value = 0.3553476858595784
$ws.Range("M8").Value = 7.323012666666667
$ws.Range("N8").Value = 21.969038
$ws.Range("O8").Value = 0.1874136762063552
$ws.Range("P8").Value = 0.1874136762063552
$ws.Range("Q8").Value = 1.251663971012
$ws.Range("R8").Value = 11.264975739108
$ws.Range("S8").Value = 0.06659701613836468
$ws.Range("T8").Value = 0.06659701613836466
# Row 9
$ws.Range("I9").Value = 0.3553476858595785
$ws.Range("J9").Value = 0.3553476858595784
$ws.Range("M9").Value = 20.55489766666667
$ws.Range("N9").Value = 61.664693
$ws.Range("O9").Value = 0.5260497436103619
$ws.Range("P9").Value = 0.5260497436103618
$ws.Range("Q9").Value = 3.513284218982
$ws.Range("R9").Value = 31.619557970838
$ws.Range("S9").Value = 0.1869305590389667
$ws.Range("T9").Value = 0.1869305590389666
# Row 10
$ws.Range("G10").Value = 0.1002363333333333
$ws.Range("H10").Value = 0.300709
$ws.Range("I10").Value = 0.2083918342229165
$ws.Range("J10").Value = 0.2083918342229164
$ws.Range("M10").Value = 5.560959666666666
$ws.Range("N10").Value = 16.682879
$ws.Range("O10").Value = 0.1423184612405788
$ws.Range("P10").Value = 0.1423184612405788
$ws.Range("Q10").Value = 0.5574102068012222
$ws.Range("R10").Value = 5.016691861211
$ws.Range("S10").Value = 0.02965800518170726
$ws.Range("T10").Value = 0.02965800518170725
# Row 11
$ws.Range("G11").Value = 0.1002363333333333
$ws.Range("H11").Value = 0.300709
$ws.Range("I11").Value = 0.2083918342229165
$ws.Range("J11").Value = 0.2083918342229164
$ws.Range("M11").Value = 5.635187000000001
$ws.Range("O11").Value = 0.1442181189427041
$ws.Range("P11").Value = 0.1442181189427041
$ws.Range("Q11").Value = 0.5648504825276668
$ws.Range("R11").Value = 5.083654342749001
$ws.Range("S11").Value = 0.03005387833464885
$ws.Range("T11").Value = 0.03005387833464884
# Row 12
$ws.Range("G12").Value = 0.1002363333333333
$ws.Range("H12").Value = 0.300709
$ws.Range("I12").Value = 0.2083918342229165
$ws.Range("J12").Value = 0.2083918342229164
$ws.Range("M12").Value = 7.323012666666667
$ws.Range("N12").Value = 21.969038
$ws.Range("O12").Value = 0.1874136762063552
$ws.Range("P12").Value = 0.1874136762063552
$ws.Range("Q12").Value = 0.7340319386602222
$ws.Range("R12").Value = 6.606287447942001
$ws.Range("S12").Value = 0.03905547974310212
$ws.Range("T12").Value = 0.03905547974310211
# Row 13
$ws.Range("G13").Value = 0.1002363333333333
$ws.Range("H13").Value = 0.300709
$ws.Range("I13").Value = 0.2083918342229165
$ws.Range("J13").Value = 0.2083918342229164
$ws.Range("M13").Value = 20.55489766666667
$ws.Range("N13").Value = 61.664693
$ws.Range("O13").Value = 0.5260497436103619
$ws.Range("P13").Value = 0.5260497436103618
$ws.Range("Q13").Value = 2.060347574148556
$ws.Range("R13").Value = 18.543128167337
$ws.Range("S13").Value = 0.1096244709634583
$ws.Range("T13").Value = 0.1096244709634582
# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.02668566666666667
$ws.Range("H14").Value = 0.080057
$ws.Range("I14").Value = 0.05547963337440523
$ws.Range("J14").Value = 0.05547963337440522
$ws.Range("M14").Value = 5.560959666666666
$ws.Range("N14").Value = 16.682879
$ws.Range("O14").Value = 0.1423184612405788
$ws.Range("P14").Value = 0.1423184612405788
$ws.Range("Q14").Value = 0.1483979160114444
$ws.Range("R14").Value = 1.335581244103
$ws.Range("S14").Value = 0.007895776052036814
$ws.Range("T14").Value = 0.007895776052036811
# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.02668566666666667
$ws.Range("H15").Value = 0.080057
$ws.Range("I15").Value = 0.05547963337440523
$ws.Range("J15").Value = 0.05547963337440522
$ws.Range("M15").Value = 5.635187000000001
$ws.Range("O15").Value = 0.1442181189427041
$ws.Range("P15").Value = 0.1442181189427041
$ws.Range("Q15").Value = 0.1503787218863334
$ws.Range("R15").Value = 1.353408496977
$ws.Range("S15").Value = 0.008001168364887592
$ws.Range("T15").Value = 0.00800116836488759
# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.02668566666666667
$ws.Range("H16").Value = 0.080057
$ws.Range("I16").Value = 0.05547963337440523
$ws.Range("J16").Value = 0.05547963337440522
$ws.Range("M16").Value = 7.323012666666667
$ws.Range("N16").Value = 21.969038
$ws.Range("O16").Value = 0.1874136762063552
$ws.Range("P16").Value = 0.1874136762063552
$ws.Range("Q16").Value = 0.1954194750184444
$ws.Range("R16").Value = 1.758775275166
$ws.Range("S16").Value = 0.01039764204527808
$ws.Range("T16").Value = 0.01039764204527808
# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.02668566666666667
$ws.Range("H17").Value = 0.080057
$ws.Range("I17").Value = 0.05547963337440523
$ws.Range("J17").Value = 0.05547963337440522
$ws.Range("M17").Value = 20.55489766666667
$ws.Range("N17").Value = 61.664693
$ws.Range("O17").Value = 0.5260497436103619
$ws.Range("P17").Value = 0.5260497436103618
$ws.Range("Q17").Value = 0.5485211475001111
$ws.Range("R17").Value = 4.936690327501
$ws.Range("S17").Value = 0.02918504691220275
$ws.Range("T17").Value = 0.02918504691220274
